$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 446, shifting existing rows 446:529 down to 447:530.
$ws.Rows.Item(446).Insert()

# Populate the newly inserted row 446 with the new record.
$ws.Cells.Item(446, 1).Value = 5
$ws.Cells.Item(446, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(446, 3).Value = "Maule"
$ws.Cells.Item(446, 4).Value = 45015
$ws.Cells.Item(446, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(446, 5).Value = 7
$ws.Cells.Item(446, 6).Value = 100112032
$ws.Cells.Item(446, 7).Value = "Zapallo italiano"
$ws.Cells.Item(446, 8).Value = "Sin especificar"
$ws.Cells.Item(446, 9).Value = "Primera"
$ws.Cells.Item(446, 10).Value = 400
$ws.Cells.Item(446, 11).Value = 5000
$ws.Cells.Item(446, 12).Value = 5000
$ws.Cells.Item(446, 13).Value = 5000
$ws.Cells.Item(446, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(446, 15).Value = "Región del Maule"
$ws.Cells.Item(446, 16).Value = 100
$ws.Cells.Item(446, 17).Value = 50
$ws.Cells.Item(446, 18).Value = "Hortaliza"
